$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table runs from row 2 through row 55 (one row per month).
# Add a new row 56 for the next month (2024-07-01), continuing the series.
# Copy the formatting of the preceding row (row 55) down into row 56 first,
# so the new row picks up the same cell style (in particular, the date
# number format used in column A), then overwrite the values.
$srcRow = $ws.Range("A55:F55")
$dstRow = $ws.Range("A56:F56")
$srcRow.Copy($dstRow)

$ws.Range("A56").Value = 45474
$ws.Range("B56").Value = -0.376
$ws.Range("C56").Value = 0.41
$ws.Range("D56").Value = -0.586
$ws.Range("E56").Value = 0.401
$ws.Range("F56").Value = 1.698
